$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header strings in row 1: "<Name>_old" -> "<Name>_FV2310",
#    "<Name>_new" -> "<Name>_FV2404". The "diff" header (column K / 11th)
#    is left untouched.
$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2) Turn the header + data range into an Excel Table (ListObject) so the
#    renamed headers are also used as the table's column names.
$rng = $ws.Range("A1:U74")
$lo = $ws.ListObjects.Add(1, $rng, [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"

# 3) Freeze the header row (split after row 1): select the first cell of
#    the scrollable area, then turn on FreezePanes - mirrors the usual
#    "View > Freeze Panes" workflow and yields state="frozen" (as opposed
#    to "frozenSplit", which is what a raw SplitRow assignment produces).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
